$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 120, pushing the existing rows 120-201 down to 122-203.
$ws.Rows.Item(120).Insert()
$ws.Rows.Item(120).Insert()

# New row 120: Kurakata / Especial
$ws.Cells.Item(120, 1).Value = 7
$ws.Cells.Item(120, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(120, 3).Value = "Ñuble"
$ws.Cells.Item(120, 4).Value = 44567
$ws.Cells.Item(120, 5).Value = 16
$ws.Cells.Item(120, 6).Value = "Fruta"
$ws.Cells.Item(120, 7).Value = 100103
$ws.Cells.Item(120, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(120, 9).Value = 100103004
$ws.Cells.Item(120, 10).Value = "Durazno"
$ws.Cells.Item(120, 11).Value = "Kurakata"
$ws.Cells.Item(120, 12).Value = "Especial"
$ws.Cells.Item(120, 13).Value = 160
$ws.Cells.Item(120, 14).Value = 14000
$ws.Cells.Item(120, 15).Value = 15000
$ws.Cells.Item(120, 16).Value = 14500
$ws.Cells.Item(120, 17).Value = "`$/caja 16 kilos empedrada"
$ws.Cells.Item(120, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(120, 19).Value = 906
$ws.Cells.Item(120, 20).Value = 16

# New row 121: Kurakata / Primera
$ws.Cells.Item(121, 1).Value = 7
$ws.Cells.Item(121, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(121, 3).Value = "Ñuble"
$ws.Cells.Item(121, 4).Value = 44567
$ws.Cells.Item(121, 5).Value = 16
$ws.Cells.Item(121, 6).Value = "Fruta"
$ws.Cells.Item(121, 7).Value = 100103
$ws.Cells.Item(121, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(121, 9).Value = 100103004
$ws.Cells.Item(121, 10).Value = "Durazno"
$ws.Cells.Item(121, 11).Value = "Kurakata"
$ws.Cells.Item(121, 12).Value = "Primera"
$ws.Cells.Item(121, 13).Value = 240
$ws.Cells.Item(121, 14).Value = 12500
$ws.Cells.Item(121, 15).Value = 13000
$ws.Cells.Item(121, 16).Value = 12750
$ws.Cells.Item(121, 17).Value = "`$/caja 16 kilos empedrada"
$ws.Cells.Item(121, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(121, 19).Value = 797
$ws.Cells.Item(121, 20).Value = 16
